$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("A1").Borders.Item(7).LineStyle = 1
$ws.Range("A1").Borders.Item(7).Weight = 4
